# Refactorizando los dirs y algunos flujos actualizados
# Remove the duplicate "Aula 2" rows for Curso D and Curso I from the CRA courses sheet.
# Row 18 = CRA_CURSOD_AULA2 (id 10846)
# Row 30 = CRA_CURSOI_AULA2 (id 10849)
# Deleting row 30 first keeps row 18's index valid for the second delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(30).Delete()
$ws.Rows.Item(18).Delete()
